# Generate Report for Handback
#
# For each localized-language sheet (zh-cn, de-de) the handback pass now
# records, per source file row:
#   - Latest Target File (col J): the handed-back markdown file, linked back
#     to its GitHub blob (same URL already used by the column-A hyperlink).
#   - Latest Handback File (col K): the generated handback xliff name.
# On the de-de sheet it also stamps "Latest Handback DateTime" (col L) with
# the handback timestamp. The Overview sheet's status column flips from
# "Ready for handoff" to "Handed back: in sync with en-US".

$wb = $excel.ActiveWorkbook

$urlMd1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3667bce13c7867b47c28a45c8d91a9d3fc024e87/e2e/4f9cad43-9d09-4a8f-8b78-e55260c6d537.md"
$urlMd2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3667bce13c7867b47c28a45c8d91a9d3fc024e87/e2e/716f89b0-12ce-4dbb-882d-7f62e50ce2b3.md"
$nameMd1 = "4f9cad43-9d09-4a8f-8b78-e55260c6d537.md"
$nameMd2 = "716f89b0-12ce-4dbb-882d-7f62e50ce2b3.md"

# ---------------------------------------------------------------------
# Status moves from "Ready for handoff" to "Handed back: in sync with
# en-US" everywhere it is shown: the Overview roll-up columns AND the
# per-file Status column on each language sheet (they share the same
# underlying text).
# ---------------------------------------------------------------------
$statusText = "Handed back: in sync with en-US"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText
$wsOverview.Range("E3").Value = $statusText
$wsOverview.Range("F3").Value = $statusText
$wsOverview.Columns.Item(5).ColumnWidth = 29.1
$wsOverview.Columns.Item(6).ColumnWidth = 29.1

# ---------------------------------------------------------------------
# zh-cn: handback xliffs generated 2017-02-09 10:11:49.
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = $statusText
$wsZh.Range("C3").Value = $statusText

$wsZh.Range("J2").Value = $nameMd1
$wsZh.Hyperlinks.Add($wsZh.Range("J2"), $urlMd1, "", "", $nameMd1)
$wsZh.Range("K2").Value = "4f9cad43-9d09-4a8f-8b78-e55260c6d537.a32078622b89f3c625e0b3170e3dac76cf31d122.zh-cn.xlf"
$wsZh.Range("L2").Value = "2017-02-09 10:11:49"

$wsZh.Range("J3").Value = $nameMd2
$wsZh.Hyperlinks.Add($wsZh.Range("J3"), $urlMd2, "", "", $nameMd2)
$wsZh.Range("K3").Value = "716f89b0-12ce-4dbb-882d-7f62e50ce2b3.0b4f02966a0683759061e45f236ac0c02dbe8d98.zh-cn.xlf"
$wsZh.Range("L3").Value = "2017-02-09 10:11:49"

$wsZh.Columns.Item(3).ColumnWidth = 29.1
$wsZh.Columns.Item(10).ColumnWidth = 39.17
$wsZh.Columns.Item(11).ColumnWidth = 39.17

# ---------------------------------------------------------------------
# de-de: handback xliffs generated slightly later, 2017-02-09 10:12:19.
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = $statusText
$wsDe.Range("C3").Value = $statusText

$wsDe.Range("J2").Value = $nameMd1
$wsDe.Hyperlinks.Add($wsDe.Range("J2"), $urlMd1, "", "", $nameMd1)
$wsDe.Range("K2").Value = "4f9cad43-9d09-4a8f-8b78-e55260c6d537.a32078622b89f3c625e0b3170e3dac76cf31d122.de-de.xlf"
$wsDe.Range("L2").Value = "2017-02-09 10:12:19"

$wsDe.Range("J3").Value = $nameMd2
$wsDe.Hyperlinks.Add($wsDe.Range("J3"), $urlMd2, "", "", $nameMd2)
$wsDe.Range("K3").Value = "716f89b0-12ce-4dbb-882d-7f62e50ce2b3.0b4f02966a0683759061e45f236ac0c02dbe8d98.de-de.xlf"
$wsDe.Range("L3").Value = "2017-02-09 10:12:19"

$wsDe.Columns.Item(3).ColumnWidth = 29.1
$wsDe.Columns.Item(10).ColumnWidth = 39.17
$wsDe.Columns.Item(11).ColumnWidth = 39.17
